$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.578048586845398
$ws.Range("B1").Value = 1.440586447715759
$ws.Range("C1").Value = 5.569024085998535
$ws.Range("D1").Value = 2.111257076263428
$ws.Range("E1").Value = 0.8993200063705444
